$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "duplicate_image_filename" column (E) header already exists in row 1.
# Populate "NA" for every data row under it (rows 2 through 21).
$ws.Range("E2:E21").Value = "NA"
